$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing the "Feeling Trendy?" question (row 2: id=0, body="Feeling Trendy?", type="bool")
$ws.Rows.Item(2).Delete()

# Update selection to match the post-deletion state (full row 2 selected, active cell A2)
$ws.Range("A2:XFD2").Select()

# Reposition the window on the landing page
$wb.Windows.Item(1).Left = 4000
